# TSTool-Command-List.xlsx update
# - Add SetPropertyFromEnsemble command (and related updates from SNODAS work)
# - Mark several existing commands as updated (Y) in columns C/D
# - Insert a new "Total commands updated" / "Remaining to do" summary section
# - Update selection/scroll state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark additional commands (columns C/D) as updated ("Y") ---
# Row 33: CreateEnsembleFromOneTimeSeries
$ws.Range("C33").Value = "Y"
$ws.Range("D33").Value = "Y"

# Row 65
$ws.Range("C65").Value = "Y"
$ws.Range("D65").Value = "Y"

# Row 68
$ws.Range("C68").Value = "Y"
$ws.Range("D68").Value = "Y"

# Row 69
$ws.Range("C69").Value = "Y"
$ws.Range("D69").Value = "Y"

# Row 117
$ws.Range("C117").Value = "Y"
$ws.Range("D117").Value = "Y"

# --- Insert a new row above the old "totals" row (old row 241) so the
#     summary block grows by one row and gains a "Total commands updated" /
#     "Remaining to do" labeled pair of rows. ---
$ws.Rows("241").Insert()

# New row 242 (was 241): totals row, now also labeled in column A.
$ws.Range("A242").Value = "Total commands updated"
$ws.Range("C242").Formula = "=COUNTIF(C2:C240,""=Y"")+COUNTIF(C2:C240,""=NA"")+COUNTIF(C2:C240,""=Z"")"
$ws.Range("D242").Formula = "=COUNTIF(D2:D240,""=Y"")+COUNTIF(D2:D240,""=NA"")+COUNTIF(D2:D240,""=Z"")"
$ws.Range("E242").Formula = "=COUNTIF(E2:E240,""=Y"")+COUNTIF(E2:E240,""=NA"")+COUNTIF(E2:E240,""=Z"")"

# New row 243 (was 242): legend row "Y = yes updated", now also labeled with
#     "Remaining to do" and new C/D formulas counting what's left.
$ws.Range("A243").Value = "Remaining to do"
$ws.Range("C243").Formula = "=ROW() - 3 - C242"
$ws.Range("D243").Formula = "=ROW() - 3 - D242"

# Rows 244/245 (were 243/244) keep their original content; the insert above
# has already shifted them down automatically.

# --- Selection / view state ---
$ws.Range("E18").Select()
